# Update time_taken values on the "data" sheet (column F, rows 2-14)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("F2").Value  = "2021-10-05 14:33:42.488474"
$ws.Range("F3").Value  = "2021-10-05 14:33:42.488482"
$ws.Range("F4").Value  = "2021-10-05 14:33:42.488485"
$ws.Range("F5").Value  = "2021-10-05 14:33:42.488487"
$ws.Range("F6").Value  = "2021-10-05 14:33:42.488490"
$ws.Range("F7").Value  = "2021-10-05 14:33:42.488493"
$ws.Range("F8").Value  = "2021-10-05 14:33:42.488496"
$ws.Range("F9").Value  = "2021-10-05 14:33:42.488498"
$ws.Range("F10").Value = "2021-10-05 14:33:42.488501"
$ws.Range("F11").Value = "2021-10-05 14:33:42.488503"
$ws.Range("F12").Value = "2021-10-05 14:33:42.488506"
$ws.Range("F13").Value = "2021-10-05 14:33:42.488508"
$ws.Range("F14").Value = "2021-10-05 14:33:42.488510"

# Add a new "metadata" worksheet after the "data" sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$meta = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Desmosomal disorders"
$meta.Range("C2").Value = 97

# data_version "0.8" must be stored as text, not a number.
# Using a leading apostrophe forces text entry, then resetting the
# Style back to "Normal" drops the quote-prefix style Excel applies,
# leaving a plain text cell with no explicit style index.
$meta.Range("D2").Value = "'0.8"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2021-07-28T09:22:41.912655Z"
$meta.Range("F2").Value = "2021-10-05 14:33:42.484624"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/97/?format=json"

# Apply the bold/border/center header style ("s=1" in the data sheet) to
# the metadata header row (and the A2 index cell), matching the style used
# for headers/index cells elsewhere in the workbook.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
